$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("A2").Value = 112234730
$ws.Range("B2").Value = 96735

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "60"

$ws.Range("K2").Value = "fullt utvecklade blad"

$ws.Range("P2").Value = "Djupsund ONO 345 m, Ög"
$ws.Range("Q2").Value = 567854
$ws.Range("R2").Value = 6511993

$ws.Range("AC2").Value = "Mer än 60 ex."

# Row 3 changes
$ws.Range("A3").Value = 112234707
$ws.Range("B3").Value = 96735

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "7"

$ws.Range("K3").Value = "överblommad"

$ws.Range("P3").Value = "Djupsund ONO 338 m, Ög"
$ws.Range("Q3").Value = 567849
$ws.Range("R3").Value = 6511980

$ws.Range("AC3").ClearContents()

# Row 4 changes
$ws.Range("B4").Value = 98980
